$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Clear rows 12-14 (A:C) - "Benachrichtigungen ansehen/löschen" + "Stornien/innerhalb 14 Tagen"
# rows entirely, keeping the existing (centered) formatting on B/C.
$ws.Range("A12:C14").ClearContents()

# Row 22: "Artikel löschen" row becomes "Order ansehen", with B/C cleared
$ws.Range("A22").Value = "Order ansehen"
$ws.Range("B22:C22").ClearContents()

# Clear rows 27-28 (A:C) - "Benachrichtigungen löschen/ansehen"
$ws.Range("A27:C28").ClearContents()

# Update the view: selection moves to A27 (matches the saved sheetView state)
$ws.Range("A27").Select()
